# Juno: check in to OLPRODLOC.
# Translate the "Charger sales report" worksheet to Spanish and rename the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "Informe de ventas"

# Translate header row (row 1).
$ws.Range("A1").Value = "Año-Trimestre"
$ws.Range("B1").Value = "Medio oeste"
$ws.Range("C1").Value = "Montaña"
$ws.Range("F1").Value = "Sudeste"
$ws.Range("G1").Value = "Oeste"

# Translate the Year-Quarter labels (column A, rows 2-9) from "Q" to "T" (Trimestre).
$ws.Range("A2").Value = "2022-T1"
$ws.Range("A3").Value = "2022-T2"
$ws.Range("A4").Value = "2022-T3"
$ws.Range("A5").Value = "2022-T4"
$ws.Range("A6").Value = "2023-T1"
$ws.Range("A7").Value = "2023-T2"
$ws.Range("A8").Value = "2023-T3"
$ws.Range("A9").Value = "2023-T4"
